# Add links to the readings listed in the Schedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedule")

# Set values in the same order the shared-string table records them in the
# target workbook (new unique strings are appended in first-use order).
$ws.Range("E3").Value  = "IMStat Chapters [2](https://openintro-ims.netlify.app/data-design) and [3](https://openintro-ims.netlify.app/data-applications)"
$ws.Range("E2").Value  = "IMStat [Chapter 1](https://openintro-ims.netlify.app/data-hello)"
$ws.Range("E5").Value  = "IMStat [Chapter 7](https://openintro-ims.netlify.app/model-slr)"
$ws.Range("E6").Value  = "IMStat [Chapter 8](https://openintro-ims.netlify.app/model-mlr)"
$ws.Range("E7").Value  = "IMStat [Chapter 9](https://openintro-ims.netlify.app/model-logistic)"
$ws.Range("E10").Value = "IMStat Chapters [16](https://openintro-ims.netlify.app/inference-one-prop), [17](https://openintro-ims.netlify.app/inference-two-props), and [18](https://openintro-ims.netlify.app/inference-tables)"
$ws.Range("E11").Value = "IMStat Chapters [19](https://openintro-ims.netlify.app/inference-one-mean), [20](https://openintro-ims.netlify.app/inference-two-means), and [21](https://openintro-ims.netlify.app/inference-paired-means)"
$ws.Range("E12").Value = "IMStat [Chapter 22](https://openintro-ims.netlify.app/inference-many-means)"
$ws.Range("E4").Value  = "R4DS [Chapter 3](https://r4ds.had.co.nz/data-visualisation.html)"

$ws.Activate()
$ws.Range("E4").Select()
